# Auto-generated edit script: updates crypto price/volume table per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.729.63"
Set-TextValue "E2" "  -3.94%  "
Set-TextValue "D3" "1.818.45"
Set-TextValue "E3" "  -2.72%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "278.90"
Set-TextValue "E5" "  -7.26%  "
Set-TextValue "D7" "0.5086"
Set-TextValue "E7" "  -4.71%  "
Set-TextValue "D8" "0.3539"
Set-TextValue "E8" "  -5.48%  "
Set-TextValue "D9" "44.39"
Set-TextValue "E9" "  -2.33%  "
Set-TextValue "D10" "0.06687"
Set-TextValue "E10" "  -6.93%  "
Set-TextValue "E11" "  -7.84%  "
Set-TextValue "D12" "0.8264"
Set-TextValue "E12" "  -7.05%  "
Set-TextValue "D13" "0.07870"
Set-TextValue "E13" "  -3.48%  "
Set-TextValue "D14" "1.816.86"
Set-TextValue "E14" "  -2.85%  "
Set-TextValue "D15" "5.081"
Set-TextValue "E15" "  -4.11%  "
Set-TextValue "D16" "87.73"
Set-TextValue "E16" "  -5.55%  "
Set-TextValue "E17" "  +0.06%  "
Set-TextValue "D18" "14.09"
Set-TextValue "E18" "  -5.11%  "
Set-TextValue "D19" "0.000008047"
Set-TextValue "E19" "  -5.30%  "
Set-TextValue "E20" "  +0.07%  "
Set-TextValue "D21" "25.771.85"
Set-TextValue "E21" "  -3.88%  "
Set-TextValue "D22" "4.754"
Set-TextValue "E22" "  -4.63%  "
Set-TextValue "D23" "10.00"
Set-TextValue "D24" "6.108"
Set-TextValue "E24" "  -4.35%  "
Set-TextValue "D25" "2.243"
Set-TextValue "E25" "  -2.98%  "
Set-TextValue "D26" "142.36"
Set-TextValue "E26" "  -2.51%  "
Set-TextValue "D27" "1.669"
Set-TextValue "E27" "  -3.76%  "
Set-TextValue "D28" "17.15"
Set-TextValue "E28" "  -4.81%  "
Set-TextValue "D29" "109.24"
Set-TextValue "E29" "  -4.18%  "
Set-TextValue "D30" "4.338"
Set-TextValue "E30" "  -8.14%  "
Set-TextValue "D31" "4.228"
Set-TextValue "E31" "  -8.65%  "
Set-TextValue "D32" "0.08761"
Set-TextValue "E32" "  -4.25%  "
Set-TextValue "D33" "0.04890"
Set-TextValue "E33" "  -2.64%  "
Set-TextValue "D34" "0.7282"
Set-TextValue "E34" "  -9.63%  "
Set-TextValue "E35" "  -3.15%  "
Set-TextValue "E36" "  -1.66%  "
Set-TextValue "D37" "3.151"
Set-TextValue "E37" "  -1.32%  "
Set-TextValue "B38" "RenderToken"
Set-TextValue "C38" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "2.361"
Set-TextValue "E38" "  -12.39%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01855"
Set-TextValue "E39" "  -5.00%  "
Set-TextValue "B40" "TheSandbox"
Set-TextValue "C40" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D40" "0.5163"
Set-TextValue "E40" "  -15.92%  "
Set-TextValue "B41" "TrustWalletToken"
Set-TextValue "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D41" "0.9718"
Set-TextValue "E41" "  -8.99%  "
Set-TextValue "B42" "Quant"
Set-TextValue "C42" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D42" "114.25"
Set-TextValue "E42" "  -0.41%  "
Set-TextValue "B43" "FraxShare"
Set-TextValue "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "6.239"
Set-TextValue "E43" "  -4.32%  "
Set-TextValue "B44" "Aptos"
Set-TextValue "C44" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D44" "8.013"
Set-TextValue "E44" "  -8.61%  "
Set-TextValue "B45" "PaxDollar"
Set-TextValue "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D45" "1.001"
Set-TextValue "E45" "  +0.09%  "
Set-TextValue "B46" "Decentraland"
Set-TextValue "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.4538"
Set-TextValue "E46" "  -12.90%  "
Set-TextValue "B47" "Algorand"
Set-TextValue "C47" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D47" "0.1369"
Set-TextValue "E47" "  -8.27%  "
Set-TextValue "B48" "Elrond"
Set-TextValue "C48" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D48" "36.47"
Set-TextValue "E48" "  -3.16%  "
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "9.200"
Set-TextValue "E49" "  -7.61%  "
Set-TextValue "B50" "NEARProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D50" "1.501"
Set-TextValue "E50" "  -9.13%  "
Set-TextValue "B51" "Cronos"
Set-TextValue "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.05839"
Set-TextValue "E51" "  -3.64%  "
